$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Neo4jData) / E (WebData) file names for rows 2-4
$ws.Range("D2").Value = "TC08_CDS_Filter_InstrumentModel-Illumina MiSeq_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC08_CDS_Filter_InstrumentModel-Illumina MiSeq_WebData.xlsx"
$ws.Range("D3").Value = "TC08_CDS_Filter_InstrumentModel-Illumina MiSeq_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC08_CDS_Filter_InstrumentModel-Illumina MiSeq_WebData.xlsx"
$ws.Range("D4").Value = "TC08_CDS_Filter_InstrumentModel-Illumina MiSeq_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC08_CDS_Filter_InstrumentModel-Illumina MiSeq_WebData.xlsx"

# Column B - per-tab query text (instrument model swapped to Illumina MiSeq)
$participantQuery = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina MiSeq']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p, s, collect(distinct samp.sample_id) as samp`nRETURN `ncoalesce(p.participant_id,'') as ``Participant ID```,`ncoalesce(s.study_name, '') as ``Study Name```,`ncoalesce(s.phs_accession,'') as ``Accession```,`ncoalesce(p.gender,'') as ``Gender```,`ncoalesce(apoc.text.join(samp, ','), '') as ``Samples```nORDER BY ``Participant ID``LIMIT 100"

$sampleQuery = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina MiSeq']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN  `n coalesce(samp.sample_id, '') as ``Sample ID```,`n coalesce(p.participant_id,'') as ``Participant ID```,`n coalesce(s.study_name, '') as ``Study Name```,`n coalesce(s.phs_accession,'') as ``Accession```,`ncoalesce(samp.sample_tumor_status,'') as ``Tumor```,`ncoalesce(samp.sample_type,'') as ``Analyte Type```nORDER By samp.sample_id LIMIT 100"

$fileQuery = "Match (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina MiSeq']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN `n    coalesce(f.file_name, '') as ``File Name```,`n    coalesce(s.study_name, '') as ``Study Name```,`n    coalesce(s.phs_accession,'') as ``Accession```,`n    coalesce(p.participant_id,'') as ``Participant ID```,`n    coalesce(samp.sample_id, '') as ``Sample ID```,`n    coalesce(f.file_type, '') as ``File Type```nORDER By f.file_name LIMIT 100"

$statsQuery = "MATCH (f:file)`nMatch (f)<--(g:genomic_info)`nWHERE g.instrument_model in ['Illumina MiSeq']`nMATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)`nWITH p,f, s, collect(distinct samp.sample_id) as samp`nRETURN`ncount(distinct s) AS Studies,`ncount(distinct p) AS Participants,`ncount(distinct samp) AS Samples,`ncount(distinct f) AS Files"

$ws.Range("B2").Value = $participantQuery
$ws.Range("B3").Value = $sampleQuery
$ws.Range("B4").Value = $fileQuery

# Column C - common stats query (same across rows 2-4)
$ws.Range("C2").Value = $statsQuery
$ws.Range("C3").Value = $statsQuery
$ws.Range("C4").Value = $statsQuery

# Update column widths (D and E) to match final bestFit widths
$ws.Columns.Item(4).ColumnWidth = 88.85546875
$ws.Columns.Item(5).ColumnWidth = 87.140625

# Update selection to D4
$ws.Range("D4").Select()
